# Layout für Map angepasst
# - Tasks sheet: mark the "Grundlegende Karte einbauen" task (row 36) as done
#   with a completion date, mirroring the other task rows above it.
# - Update the active selection to C37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Zustand (state) column -> "done", same text used throughout column C.
$ws.Range("C36").Value = "done"

# abgeschlossen am (completed on) column -> copy formatting from the row
# above (D35) so the new date cell keeps the same date number format.
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D36").Value = Get-Date -Year 2019 -Month 3 -Day 21 -Hour 0 -Minute 0 -Second 0

# Move the selection like Excel would after entering data in D36.
$ws.Range("C37").Select() | Out-Null
